$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("Tipo" currently lives there),
# shifting the existing "Tipo" column to E, to make room for a new "MAE" column.
$ws.Range("D1").EntireColumn.Insert()

# Copy the header formatting (bold, border, alignment) from the neighboring
# "R2" header cell onto the new "MAE" header cell.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats

# New header text for the inserted column
$ws.Range("D1").Value = "MAE"

# New MAE value for the data row
$ws.Range("D2").Value = 0.8429226117853587
